$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'72.351.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.69%  "
$ws.Range("D3").Value = "'2.665.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'602.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").Value = "'178.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.61%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.524"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("D9").Value = "'2.663.86"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.46%  "
$ws.Range("D10").Value = "'0.171"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.00%  "
$ws.Range("E11").Value = "  +2.36%  "
$ws.Range("D12").Value = "'0.356"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.65%  "
$ws.Range("D13").Value = "'5.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").Value = "'3.152.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.30%  "
$ws.Range("D15").Value = "'0.0000186"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.74%  "
$ws.Range("D16").Value = "'72.272.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.53%  "
$ws.Range("D17").Value = "'26.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.40%  "
$ws.Range("D18").Value = "'2.677.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.52%  "
$ws.Range("D19").Value = "'11.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.03%  "
$ws.Range("D20").Value = "'7.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").Value = "'370.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.04%  "
$ws.Range("D22").Value = "'4.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("D23").Value = "'2.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.52%  "
$ws.Range("D24").Value = "'72.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").Value = "'4.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.68%  "
$ws.Range("D27").Value = "'9.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.04%  "
$ws.Range("D28").Value = "'2.804.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").Value = "'0.0₃0943"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.38%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").Value = "'519.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.11%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'8.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("D33").Value = "'1.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("D34").Value = "'1.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "'163.27"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.34%  "
$ws.Range("D37").Value = "'19.45"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.17%  "
$ws.Range("E38").Value = "  +0.93%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.110"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.50%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "'1.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").Value = "'1.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.74%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").Value = "'5.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.57%  "
$ws.Range("D44").Value = "'2.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.45%  "
$ws.Range("D45").Value = "'0.334"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").Value = "'39.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.05%  "
$ws.Range("D47").Value = "'152.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.36%  "
$ws.Range("D48").Value = "'3.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.48%  "
$ws.Range("D49").Value = "'0.546"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.55%  "
$ws.Range("D50").Value = "'1.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.68%  "
$ws.Range("D51").Value = "'0.0767"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.49%  "
